$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'24.814.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.46%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.659.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -5.73%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +1.09%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("B5").Value = "'BNB"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'308.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.53%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("B6").Value = "'USDC"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.99%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.3647"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -5.21%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.3328"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -9.06%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'47.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -7.90%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'1.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -6.97%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07198"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -6.85%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.11%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'6.109"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -5.69%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'20.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -8.24%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'6.740"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.09%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'1.655.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.19%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.00001082"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -7.40%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.17%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.06610"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.48%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'80.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -7.49%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'16.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -6.31%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'6.063"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -6.49%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'12.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.41%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'24.771.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.27%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'2.415"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.63%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'2.624"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -11.19%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'148.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.91%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'19.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.97%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'128.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.28%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'1.839.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.13%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'1.197"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.58%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'4.134"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.43%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'6.279"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -11.07%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'1.733"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.85%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.08490"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.53%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'13.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -8.74%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'5.313"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.26%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.06313"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -6.71%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.02311"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.69%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'8.627"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -8.24%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.2125"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.88%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'1.233"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.21%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.6168"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.36%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.98%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'13.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.31%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'3.774"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.53%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.5869"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -8.04%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'2.012"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -8.06%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'124.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -6.04%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.07104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.35%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'75.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.15%  "
$ws.Range("E51").Style = "Normal"
